$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 4842
$ws.Range("I20").Value = 1410.4
$ws.Range("J20").Value = 22000
$ws.Range("K20").Value = 1410.4
$ws.Range("L20").Value = 22000
$ws.Range("M20").Value = -1180.4
$ws.Range("N20").Value = -22460
$ws.Range("H35").Value = 4842
$ws.Range("I35").Value = 1410.4
$ws.Range("J35").Value = 22000
$ws.Range("K35").Value = 1410.4
$ws.Range("L35").Value = 22000
$ws.Range("M35").Value = -1031.4
$ws.Range("N35").Value = -22758

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5084.3335
$ws.Range("I26").Value = 502.33334
$ws.Range("K26").Value = 502.33334
$ws.Range("M26").Value = -172.33334
$ws.Range("H28").Value = 6868.3076
$ws.Range("I28").Value = 4797.1816
$ws.Range("J28").Value = 18259.5
$ws.Range("K28").Value = 4797.1816
$ws.Range("L28").Value = 18259.5
$ws.Range("M28").Value = -4605.1816
$ws.Range("N28").Value = -18643.5
$ws.Range("H93").Value = 24000
$ws.Range("J93").Value = 24000
$ws.Range("L93").Value = 24000
$ws.Range("N93").Value = -28992
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21802
$ws.Range("H99").Value = 6868.3076
$ws.Range("I99").Value = 4797.1816
$ws.Range("J99").Value = 18259.5
$ws.Range("K99").Value = 4797.1816
$ws.Range("L99").Value = 18259.5
$ws.Range("M99").Value = -1802.1816
$ws.Range("N99").Value = -24249.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 22783.5
$ws.Range("J21").Value = 22783.5
$ws.Range("L21").Value = 22783.5
$ws.Range("N21").Value = -23255.5
$ws.Range("H26").Value = 24400
$ws.Range("I26").Value = 20800
$ws.Range("J26").Value = 28000
$ws.Range("K26").Value = 20800
$ws.Range("L26").Value = 28000
$ws.Range("M26").Value = -20508
$ws.Range("N26").Value = -28584
$ws.Range("H54").Value = 2582.5
$ws.Range("I54").Value = 1229.125
$ws.Range("J54").Value = 7996
$ws.Range("K54").Value = 1229.125
$ws.Range("L54").Value = 7996
$ws.Range("M54").Value = -745.125
$ws.Range("N54").Value = -8964
$ws.Range("H93").Value = 19800
$ws.Range("J93").Value = 19800
$ws.Range("L93").Value = 19800
$ws.Range("N93").Value = -23544
$ws.Range("H96").Value = 4500
$ws.Range("I96").Value = 4500
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4500
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = -1754
$ws.Range("H134").Value = 3653.3044
$ws.Range("I134").Value = 3721.8462
$ws.Range("J134").Value = 3271.4285
$ws.Range("K134").Value = 11165.5386
$ws.Range("L134").Value = 9814.2855
$ws.Range("M134").Value = -8630.5386
$ws.Range("N134").Value = -14884.2855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1110818.9
$ws.Range("I31").Value = 722.5263
$ws.Range("J31").Value = 3641838.5
$ws.Range("K31").Value = 722.5263
$ws.Range("L31").Value = 3641838.5
$ws.Range("M31").Value = -427.5263
$ws.Range("N31").Value = -3642428.5
$ws.Range("H33").Value = 2130.1428
$ws.Range("I33").Value = 1651.8334
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 1651.8334
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = -1272.8334
$ws.Range("N33").Value = -5758
$ws.Range("H34").Value = 1110818.9
$ws.Range("I34").Value = 722.5263
$ws.Range("J34").Value = 3641838.5
$ws.Range("K34").Value = 722.5263
$ws.Range("L34").Value = 3641838.5
$ws.Range("M34").Value = -520.5263
$ws.Range("N34").Value = -3642242.5
$ws.Range("H39").Value = 2500
$ws.Range("I39").Value = 2500
$ws.Range("K39").Value = 2500
$ws.Range("M39").Value = -2109
$ws.Range("H49").Value = 2500
$ws.Range("I49").Value = 2500
$ws.Range("K49").Value = 2500
$ws.Range("M49").Value = -2318
$ws.Range("H58").Value = 5269.923
$ws.Range("I58").Value = 5958.1055
$ws.Range("J58").Value = 3402
$ws.Range("K58").Value = 5958.1055
$ws.Range("L58").Value = 3402
$ws.Range("M58").Value = -5755.1055
$ws.Range("N58").Value = -3808
$ws.Range("H114").Value = 14490
$ws.Range("J114").Value = 14490
$ws.Range("L114").Value = 14490
$ws.Range("N114").Value = -23168
$ws.Range("H136").Value = 5269.923
$ws.Range("I136").Value = 5958.1055
$ws.Range("J136").Value = 3402
$ws.Range("K136").Value = 17874.3165
$ws.Range("L136").Value = 10206
$ws.Range("M136").Value = -15324.3165
$ws.Range("N136").Value = -15306

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2550.3333
$ws.Range("I7").Value = 2151
$ws.Range("K7").Value = 2151
$ws.Range("M7").Value = -2039
$ws.Range("H8").Value = 2550.3333
$ws.Range("I8").Value = 2151
$ws.Range("K8").Value = 2151
$ws.Range("M8").Value = -2012
$ws.Range("H20").Value = 9055
$ws.Range("I20").Value = 9055
$ws.Range("K20").Value = 9055
$ws.Range("M20").Value = -8810
$ws.Range("H25").Value = 35004.5
$ws.Range("J25").Value = 35004.5
$ws.Range("L25").Value = 35004.5
$ws.Range("N25").Value = -36062.5
$ws.Range("H94").Value = 24500
$ws.Range("J94").Value = 24500
$ws.Range("L94").Value = 24500
$ws.Range("N94").Value = -25852
$ws.Range("H96").Value = 20261
$ws.Range("J96").Value = 20261
$ws.Range("L96").Value = 20261
$ws.Range("N96").Value = -25753
$ws.Range("H100").Value = 34980
$ws.Range("J100").Value = 34980
$ws.Range("L100").Value = 34980
$ws.Range("N100").Value = -37144
$ws.Range("H132").Value = 2910597.8
$ws.Range("I132").Value = 4169988.5
$ws.Range("J132").Value = 4311.6924
$ws.Range("K132").Value = 12509965.5
$ws.Range("L132").Value = 12935.0772
$ws.Range("M132").Value = -12507435.5
$ws.Range("N132").Value = -17995.0772

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1950
$ws.Range("J19").Value = 1950
$ws.Range("L19").Value = 1950
$ws.Range("N19").Value = -2290

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 4708.5
$ws.Range("I61").Value = 3637.75
$ws.Range("J61").Value = 6850
$ws.Range("K61").Value = 3637.75
$ws.Range("L61").Value = 6850
$ws.Range("M61").Value = -3345.75
$ws.Range("N61").Value = -7434
$ws.Range("H102").Value = 29666.666
$ws.Range("J102").Value = 29666.666
$ws.Range("L102").Value = 29666.666
$ws.Range("N102").Value = -36156.666
$ws.Range("H104").Value = 20370
$ws.Range("J104").Value = 20370
$ws.Range("L104").Value = 20370
$ws.Range("N104").Value = -27358
$ws.Range("H106").Value = 20985
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H126").Value = 26316660
$ws.Range("I126").Value = 28572202
$ws.Range("J126").Value = 1994
$ws.Range("K126").Value = 85716606
$ws.Range("L126").Value = 5982
$ws.Range("M126").Value = -85714136
$ws.Range("N126").Value = -10922
$ws.Range("H136").Value = 16749665
$ws.Range("I136").Value = 21514076
$ws.Range("J136").Value = 5553300
$ws.Range("K136").Value = 64542228
$ws.Range("L136").Value = 16659900
$ws.Range("M136").Value = -64539678
$ws.Range("N136").Value = -16665000
